$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20 (Sofiane Belle) gets a "Taille" (height) value in column E
$ws.Range("E20").Value = "1m81"

# Row 22 (Naim Dhib) date of birth correction
$ws.Range("C22").Value = 35854

# Update the active selection to F27
$ws.Range("F27").Select()
